$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.420.93"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.82"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7068"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.22"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3139"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07856"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.48"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08024"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.903.95"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.198"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.28"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7003"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.450"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008360"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.404.35"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.31"
$ws.Range("E19").Value = "  +3.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.124.48"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.604"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1559"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.004"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.73"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.322"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.212"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05305"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.887"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7532"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.166"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01877"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267.35"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.739"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8984"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.13"
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.967"
$ws.Range("E43").Value = "  -8.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.47"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.025.72"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.548"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.789"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5169"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4307"
$ws.Range("E51").Value = "  -1.75%  "
